$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8, shifting existing row 8 (and below) down by one.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new superblock record.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "[376, 402, 373, 375]"
$ws.Cells.Item(8, 3).Value = "[397, 401, 405, 446]"
$ws.Cells.Item(8, 4).Value = "[397, 400, 401, 370, 403, 372, 404, 374, 405, 381, 446, 447]"
$ws.Cells.Item(8, 5).Value = "[397, 401, 405, 446]"
